$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Daily COVID-19 figures update (new cases / new deaths revisions) ---

# 16.01.2021 - new positive cases revised
$ws.Range("C326").Value = 59

# 01.02.2021 - new positive cases revised
$ws.Range("C342").Value = 124

# 02.02.2021 - new positive cases revised
$ws.Range("C343").Value = 103

# 03.02.2021 - new positive cases revised
$ws.Range("C344").Value = 146

# 04.02.2021 - new positive cases + extra-hospital death revised
$ws.Range("C345").Value = 59

# M345 is formatted as Text ("@"); flip to a plain numeric format so the
# value is stored as a real number instead of a text string, then restore
# the Text format (mirrors how the source data was originally entered).
$ws.Range("M345").NumberFormat = "0"
$ws.Range("M345").Value = 1
$ws.Range("M345").NumberFormat = "@"

# 05.02.2021 - brand new row of data added
$ws.Range("C346").Value = 9
$ws.Range("E346").Value = 14
$ws.Range("F346").Value = 8
$ws.Range("G346").Value = 108

$ws.Range("L346").NumberFormat = "0"
$ws.Range("L346").Value = 0
$ws.Range("L346").NumberFormat = "@"

$ws.Range("M346").NumberFormat = "0"
$ws.Range("M346").Value = 0
$ws.Range("M346").NumberFormat = "@"

# Update the active selection to reflect where the author left off editing
$ws.Range("A2").Select()
